$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 775, pushing existing rows 775-816 down to 776-817
$ws.Rows.Item(775).Insert()

# Force column A to stay as literal text "2026/02/09" (not auto-converted to a date serial)
$ws.Cells.Item(775, 1).NumberFormat = "@"
$ws.Cells.Item(775, 1).Value = "2026/02/09"
$ws.Cells.Item(775, 1).ClearFormats()

$ws.Cells.Item(775, 2).Value = "月"
$ws.Cells.Item(775, 3).Value = 6
$ws.Cells.Item(775, 4).Value = 105
